# Localization status report regeneration ("Generate Report for Archive").
#
# The CI job re-ran and the single outstanding file's status moved on from
# "Ready for handoff" to "In Translation". That status string is mirrored on
# the "Overview" sheet (once per target locale column, zh-cn / de-de) and on
# each locale's own detail sheet (its "Status" column). Because the new text
# is shorter than the old text, the locale-status columns that were
# previously sized to fit "Ready for handoff" re-fit themselves to the new,
# narrower "In Translation" text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn status is column E, de-de status is column F ---
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# --- zh-cn detail sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
$wsZhCn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# --- de-de detail sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}
$wsDeDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
